$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values with re-pulled / recalculated data.
$updates = @{
    "F2"  = -6
    "F3"  = -2
    "F4"  = 2
    "F5"  = 6
    "F6"  = -7
    "F7"  = -2
    "F8"  = 5
    "F10" = 7
    "F11" = -3
    "F12" = 8
    "F14" = 6
    "F15" = -2
    "F16" = 0
    "F17" = -5
    "F19" = -5
    "F21" = -1
    "F23" = -4
    "F24" = 1
    "F25" = -5
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
